$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geometry")

# Insert a new column before AN, shifting AN:AQ -> AO:AR
$ws.Columns("AN:AN").Insert()

# New header cell: apply same formatting as the neighboring header cells
# (bold font, thin border, centered horizontal / top vertical alignment)
$ws.Range("AN1").Value = "solidity"
$ws.Range("AN1").Font.Bold = $true
$ws.Range("AN1").Borders.LineStyle = 1
$ws.Range("AN1").HorizontalAlignment = -4108
$ws.Range("AN1").VerticalAlignment = -4160

# New data cell for row 2
$ws.Range("AN2").Value = "[1.42997704 1.70997375]"
